$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 0.6189607300944256
$ws.Range("C4").Value = 0.63
$ws.Range("D4").Value = 0.6161130996815171
$ws.Range("E4").Value = 0.6164999999999999
$ws.Range("F4").Value = 0.6013256889962927
$ws.Range("G4").Value = 0.619
$ws.Range("H4").Value = 0.595721735522868
$ws.Range("I4").Value = 0.5925
$ws.Range("J4").Value = 0.5048350243657124
$ws.Range("K4").Value = 0.4980000000000001
$ws.Range("L4").Value = 0.5172935460703159
$ws.Range("M4").Value = 0.5149999999999999
$ws.Range("N4").Value = 0.6426097101964803
$ws.Range("O4").Value = 0.6709999999999999
$ws.Range("P4").Value = 0.6231257749317758
$ws.Range("Q4").Value = 0.632
$ws.Range("R4").Value = 0.623661845275809
$ws.Range("S4").Value = 0.625
$ws.Range("T4").Value = 0.631191903743993
$ws.Range("U4").Value = 0.626
$ws.Range("V4").Value = 0.6019655784651958
$ws.Range("W4").Value = 0.621
$ws.Range("X4").Value = 0.5943432857260397
$ws.Range("Y4").Value = 0.592
$ws.Range("Z4").Value = 0.6419193021651175
$ws.Range("AA4").Value = 0.6689999999999999
$ws.Range("AB4").Value = 0.6238375418450141
$ws.Range("AC4").Value = 0.6319999999999999
$ws.Range("B5").Value = 0.6543492593164338
$ws.Range("C5").Value = 0.826
$ws.Range("D5").Value = 0.5438920968120968
$ws.Range("E5").Value = 0.5675000000000001
$ws.Range("F5").Value = 0.236562059063871
$ws.Range("G5").Value = 0.177
$ws.Range("H5").Value = 0.5437218406593407
$ws.Range("I5").Value = 0.521
$ws.Range("J5").Value = 0.6710048353787262
$ws.Range("K5").Value = 0.9240000000000002
$ws.Range("L5").Value = 0.5292454441750143
$ws.Range("M5").Value = 0.548
$ws.Range("N5").Value = 0.2804802005151907
$ws.Range("O5").Value = 0.219
$ws.Range("P5").Value = 0.5520071571248042
$ws.Range("Q5").Value = 0.5355
$ws.Range("R5").Value = 0.6550623607403533
$ws.Range("S5").Value = 0.845
$ws.Range("T5").Value = 0.5370930571893868
$ws.Range("U5").Value = 0.5580000000000001
$ws.Range("V5").Value = 0.2518515978656548
$ws.Range("W5").Value = 0.19
$ws.Range("X5").Value = 0.545344344937166
$ws.Range("Y5").Value = 0.5255
$ws.Range("Z5").Value = 0.2822188676716513
$ws.Range("AA5").Value = 0.222
$ws.Range("AB5").Value = 0.5511832757406527
$ws.Range("AC5").Value = 0.5349999999999999
$ws.Range("B6").Value = 0.6324155858997562
$ws.Range("C6").Value = 0.6060000000000001
$ws.Range("D6").Value = 0.6728368859275959
$ws.Range("E6").Value = 0.6529999999999999
$ws.Range("F6").Value = 0.6912000924665993
$ws.Range("G6").Value = 0.6830000000000001
$ws.Range("H6").Value = 0.7190885525075216
$ws.Range("I6").Value = 0.701
$ws.Range("J6").Value = 0.5035539918153584
$ws.Range("K6").Value = 0.502
$ws.Range("L6").Value = 0.5142681497684091
$ws.Range("M6").Value = 0.513
$ws.Range("N6").Value = 0.7002753360393046
$ws.Range("O6").Value = 0.6860000000000002
$ws.Range("P6").Value = 0.7449875039821274
$ws.Range("Q6").Value = 0.7190000000000001
$ws.Range("R6").Value = 0.6454031032363178
$ws.Range("S6").Value = 0.6289999999999999
$ws.Range("T6").Value = 0.6757038601739064
$ws.Range("U6").Value = 0.6609999999999999
$ws.Range("V6").Value = 0.6889474327614385
$ws.Range("W6").Value = 0.6849999999999999
$ws.Range("X6").Value = 0.7126369760128052
$ws.Range("Y6").Value = 0.6984999999999999
$ws.Range("Z6").Value = 0.7035751480782269
$ws.Range("AA6").Value = 0.6940000000000001
$ws.Range("AB6").Value = 0.7430336976205992
$ws.Range("AC6").Value = 0.72
